# Regenerate merged AHB files:
#  - rename the "_old"/"_new" column-header suffixes to "_FV2310"/"_FV2404"
#  - freeze the header row (row 1)
#  - turn the A1:U58 range into a native Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Rename header row (row 1) text -------------------------------------
# Columns A-J were "<Name>_old", K was "diff", L-U were "<Name>_new".
# They become "<Name>_FV2310", "diff", "<Name>_FV2404" respectively.
$headers = @(
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# --- 2. Freeze the header row (pane split after row 1) ----------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert A1:U58 into an Excel Table named "Table1" -------------------
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Output "Edit applied"
